$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.506.11"
$ws.Range("E2").Value = "  -4.39%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.944.45"
$ws.Range("E3").Value = "  -1.23%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.09"
$ws.Range("E5").Value = "  -4.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.69"
$ws.Range("E6").Value = "  +4.94%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("E8").Value = "  +2.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.937.47"
$ws.Range("E9").Value = "  -1.38%  "

# Row 10
$ws.Range("E10").Value = "  -4.28%  "

# Row 11
$ws.Range("E11").Value = "  -5.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  +1.35%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000221"
$ws.Range("E13").Value = "  -0.48%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.71"
$ws.Range("E14").Value = "  +0.83%  "

# Row 15
$ws.Range("E15").Value = "  +0.65%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.431.78"
$ws.Range("E16").Value = "  -1.13%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.76"
$ws.Range("E17").Value = "  +9.12%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.940.60"
$ws.Range("E18").Value = "  -1.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.506.16"
$ws.Range("E19").Value = "  -4.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "414.90"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.08"
$ws.Range("E21").Value = "  -0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("E22").Value = "  +3.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.93"
$ws.Range("E23").Value = "  -1.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.96"
$ws.Range("E24").Value = "  +2.77%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.01"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("E27").Value = "  +0.02%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.48"
$ws.Range("E28").Value = "  -2.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.50"
$ws.Range("E29").Value = "  +3.37%  "

# Row 30
$ws.Range("E30").Value = "  +4.55%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.09"
$ws.Range("E31").Value = "  -0.82%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.03"
$ws.Range("E32").Value = "  -1.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.101"
$ws.Range("E33").Value = "  +8.22%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.61"
$ws.Range("E34").Value = "  +0.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.10"
$ws.Range("E35").Value = "  -3.40%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.935"
$ws.Range("E36").Value = "  -1.48%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.49"
$ws.Range("E37").Value = "  -1.91%  "

# Row 38
$ws.Range("E38").Value = "  +3.91%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.43"
$ws.Range("E39").Value = "  +6.02%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0349"
$ws.Range("E40").Value = "  -2.99%  "

# Row 41
$ws.Range("E41").Value = "  -0.57%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.53"
$ws.Range("E42").Value = "  +3.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "377.68"
$ws.Range("E43").Value = "  -1.91%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.660.20"
$ws.Range("E44").Value = "  +1.24%  "

# Row 45
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.238"
$ws.Range("E46").Value = "  +1.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.28"
$ws.Range("E47").Value = "  +2.31%  "

# Row 48
$ws.Range("E48").Value = "  +2.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.98"
$ws.Range("E49").Value = "  -0.32%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.33"
$ws.Range("E50").Value = "  -0.27%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.00"
$ws.Range("E51").Value = "  -0.43%  "
